$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1299572217078198
$ws.Range("C2").Value = 2
$ws.Range("B3").Value = 0.1078277705427478
$ws.Range("C3").Value = 2
$ws.Range("B4").Value = 0.07221859330515201
$ws.Range("C4").Value = 2
$ws.Range("B5").Value = 0.05654472158350232
$ws.Range("C5").Value = 2
$ws.Range("B6").Value = 0.1277588633555749
$ws.Range("C6").Value = 2
$ws.Range("B7").Value = 0.1249038267153906
$ws.Range("C7").Value = 2
$ws.Range("B8").Value = 0.0663644048970903
$ws.Range("C8").Value = 2
$ws.Range("B9").Value = 0.08036746554517162
$ws.Range("C9").Value = 2
$ws.Range("B10").Value = 0.1772433055648895
$ws.Range("C10").Value = 1
$ws.Range("B11").Value = 0.1700562114568257
$ws.Range("C11").Value = 1
$ws.Range("B12").Value = 0.05831001269103774
$ws.Range("C12").Value = 2
$ws.Range("B13").Value = 0.09228085546918181
$ws.Range("C13").Value = 2
$ws.Range("B14").Value = 0.123270142683938
$ws.Range("C14").Value = 2
$ws.Range("B15").Value = 0.1652669518341842
$ws.Range("C15").Value = 1
$ws.Range("B16").Value = 0.1748656137820481
$ws.Range("C16").Value = 1
$ws.Range("B17").Value = 0.1083166556028633
$ws.Range("C17").Value = 2
$ws.Range("B18").Value = 0.07578099330335455
$ws.Range("C18").Value = 2
$ws.Range("B19").Value = 0.06753406744978903
$ws.Range("C19").Value = 2
$ws.Range("B20").Value = 0.1228717730722009
$ws.Range("C20").Value = 2
$ws.Range("B21").Value = 0.1140943065109784
$ws.Range("C21").Value = 2
$ws.Range("B22").Value = 0.1126651674067496
$ws.Range("C22").Value = 2
$ws.Range("B23").Value = 0.1110862561259321
$ws.Range("C23").Value = 2
$ws.Range("B24").Value = 0.1240096001905824
$ws.Range("C24").Value = 2
$ws.Range("B25").Value = 0.1754319859606615
$ws.Range("C25").Value = 1
$ws.Range("B26").Value = 0.1043092275492984
$ws.Range("C26").Value = 2
$ws.Range("B27").Value = 0.09192056229354915
$ws.Range("C27").Value = 2
$ws.Range("B28").Value = 0.1717763287544199
$ws.Range("C28").Value = 1
$ws.Range("B29").Value = 0.08747691312893863
$ws.Range("C29").Value = 2
$ws.Range("B30").Value = 0.06450967741542374
$ws.Range("C30").Value = 2
$ws.Range("B31").Value = 0.06985293737744432
$ws.Range("C31").Value = 2
$ws.Range("B32").Value = 0.07703882408952155
$ws.Range("C32").Value = 2
$ws.Range("B33").Value = 0.2656213135465063
$ws.Range("C33").Value = 1
$ws.Range("B34").Value = 0.1786932185669706
$ws.Range("C34").Value = 1
$ws.Range("B35").Value = 0.1369599947064933
$ws.Range("C35").Value = 2
$ws.Range("B36").Value = 0.105364705595443
$ws.Range("C36").Value = 2
$ws.Range("B37").Value = 0.1238043828202657
$ws.Range("C37").Value = 2
$ws.Range("B38").Value = 0.1801552205123295
$ws.Range("C38").Value = 1
$ws.Range("B39").Value = 0.09515722443765431
$ws.Range("C39").Value = 2
$ws.Range("B40").Value = 0.2056898789871626
$ws.Range("C40").Value = 1
$ws.Range("B41").Value = 0.2099854374346849
$ws.Range("C41").Value = 1
$ws.Range("B42").Value = 0.07877332615291859
$ws.Range("C42").Value = 2
$ws.Range("B43").Value = 0.1086876268939417
$ws.Range("C43").Value = 2
$ws.Range("B44").Value = 0.05569171257079083
$ws.Range("C44").Value = 2
$ws.Range("B45").Value = 0.2406788059593526
$ws.Range("C45").Value = 1
$ws.Range("B46").Value = 0.107352402364013
$ws.Range("C46").Value = 2
$ws.Range("B47").Value = 0.078254023556966
$ws.Range("C47").Value = 2
$ws.Range("B48").Value = 0.1805443419610509
$ws.Range("C48").Value = 1
$ws.Range("B49").Value = 0.2123287088113983
$ws.Range("C49").Value = 1
$ws.Range("B50").Value = 0.1383173930869928
$ws.Range("C50").Value = 2
$ws.Range("B51").Value = 0.08411576694153078
$ws.Range("C51").Value = 2
$ws.Range("B52").Value = 0.04064479191687972
$ws.Range("C52").Value = 2
$ws.Range("B53").Value = 0.08391397672689525
$ws.Range("C53").Value = 2
